# Updated code, steps + uploaded video for smart irrigation (moisture)

$wb = $excel.ActiveWorkbook

# --- Un-minimize the workbook window ---
$wb.Windows.Item(1).WindowState = -4143  # xlNormal

# --- Sheet: Overview ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("B5").Value = "5//5"
$ws1.Range("B6").Value = "6//6"
$ws1.Range("B6").NumberFormat = "d-mmm"
$ws1.Range("B7").Value = "Updated"

# --- Sheet: Components ---
$ws2 = $wb.Worksheets.Item("Components")
$ws2.Range("B6").Interior.Color = 5296274   # RGB(146,208,80)
[void]$ws2.Range("B6").Select()

# --- Sheet: Steps ---
$ws3 = $wb.Worksheets.Item("Steps")
$ws3.Range("A15").Interior.Color = 5287936  # RGB(0,176,80)
$ws3.Range("A16").Interior.Color = 5287936  # RGB(0,176,80)
$ws3.Range("B16").Value = "Upload code"
$ws3.Range("C15").Value = "Not sure how to frame steps?"
$ws3.Range("C15").Font.Bold = $true
$ws3.Range("C15").Interior.Color = 65535    # RGB(255,255,0)
[void]$ws3.Range("C15").Select()

# --- Restore the originally active sheet/tab ---
[void]$ws1.Activate()
